# Weekly update: insert a new data row for the latest week (2023-10-04)
# into the "Haba" sheet, pushing the existing rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 26; this shifts rows 26..66
# down to 27..67 (matching the rest of the diff) and grows the
# dimension from A1:R66 to A1:R67 automatically.
$ws.Rows.Item(26).Insert()

# Populate the freshly inserted row 26 with the new week's record.
$ws.Cells.Item(26, 1).Value  = 11
$ws.Cells.Item(26, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(26, 3).Value  = "Bíobío"
$ws.Cells.Item(26, 4).Value  = 45203
$ws.Cells.Item(26, 5).Value  = 8
$ws.Cells.Item(26, 6).Value  = 100112026
$ws.Cells.Item(26, 7).Value  = "Haba"
$ws.Cells.Item(26, 8).Value  = "Sin especificar"
$ws.Cells.Item(26, 9).Value  = "Primera"
$ws.Cells.Item(26, 10).Value = 100
$ws.Cells.Item(26, 11).Value = 11000
$ws.Cells.Item(26, 12).Value = 12000
$ws.Cells.Item(26, 13).Value = 11500
$ws.Cells.Item(26, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(26, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(26, 16).Value = 460
$ws.Cells.Item(26, 17).Value = 25
$ws.Cells.Item(26, 18).Value = "Hortaliza"

# Keep the date column formatted the same way as the rest of column D.
$ws.Cells.Item(26, 4).NumberFormat = $ws.Cells.Item(27, 4).NumberFormat
